# Sync file from Google Drive
# Updates EstimatedTimeOfArrival (F), Load (I), TypeOfBus (L) and Monitored (J)
# values on the NextBus2 and NextBus3 sheets to match the refreshed export.

$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("NextBus2")
$ws3 = $wb.Worksheets.Item("NextBus3")

# --- NextBus2 ---
$ws2.Range("F2").Value = 45684.40355324074
$ws2.Range("I2").Value = "SDA"

$ws2.Range("F3").Value = 45684.40421296296

$ws2.Range("F4").Value = 45684.39903935185

$ws2.Range("F5").Value = 45684.39681712963

$ws2.Range("F6").Value = 45684.39819444445
$ws2.Range("L6").Value = "DD"

$ws2.Range("F7").Value = 45684.40039351852
$ws2.Range("I7").Value = "SDA"

$ws2.Range("F8").Value = 45684.3962962963

$ws2.Range("F9").Value = 45684.39949074074
$ws2.Range("I9").Value = "SDA"

$ws2.Range("F10").Value = 45684.40387731481

$ws2.Range("F11").Value = 45684.39221064815
$ws2.Range("L11").Value = "SD"

$ws2.Range("F12").Value = 45684.40273148148
$ws2.Range("J12").Value = 0

$ws2.Range("F13").Value = 45684.40793981482

$ws2.Range("F14").Value = 45684.3934375

$ws2.Range("F15").Value = 45684.40432870371

# --- NextBus3 ---
$ws3.Range("F2").Value = 45684.41082175926

$ws3.Range("F3").Value = 45684.40878472223
$ws3.Range("L3").Value = "BD"

$ws3.Range("F4").Value = 45684.4047337963

$ws3.Range("F5").Value = 45684.40400462963

$ws3.Range("F6").Value = 45684.40660879629
$ws3.Range("L6").Value = "SD"

$ws3.Range("F7").Value = 45684.40424768518
$ws3.Range("L7").Value = "DD"

$ws3.Range("F8").Value = 45684.39987268519
$ws3.Range("I8").Value = "SEA"
$ws3.Range("L8").Value = "DD"

$ws3.Range("F9").Value = 45684.40621527778

$ws3.Range("F10").Value = 45684.41209490741

$ws3.Range("F11").Value = 45684.39475694444

$ws3.Range("F12").Value = 45684.41314814815

$ws3.Range("F13").Value = 45684.41668981482

$ws3.Range("F14").Value = 45684.39988425926

$ws3.Range("F15").Value = 45684.41159722222
$ws3.Range("J15").Value = 0
